$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2
$ws.Cells.Item(2, 3).Value = 0.5245901639344263
$ws.Cells.Item(2, 10).Value = 0.01639344262295082
$ws.Cells.Item(2, 16).Value = 0.1278688524590164
$ws.Cells.Item(2, 19).Value = 0.1311475409836066
$ws.Cells.Item(3, 2).Value = 0.01219512195121951
$ws.Cells.Item(3, 3).Value = 0.02439024390243903
$ws.Cells.Item(3, 10).Value = 0.006097560975609756
$ws.Cells.Item(3, 16).Value = 0.725609756097561
$ws.Cells.Item(3, 19).Value = 0.2317073170731707
$ws.Cells.Item(4, 16).Value = 0.7804878048780488
$ws.Cells.Item(4, 19).Value = 0.2195121951219512
$ws.Cells.Item(6, 2).Value = 0.0625
$ws.Cells.Item(6, 4).Value = 0.01442307692307692
$ws.Cells.Item(6, 6).Value = 0.02403846153846154
$ws.Cells.Item(6, 10).Value = 0.2451923076923077
$ws.Cells.Item(6, 15).Value = 0.02403846153846154
$ws.Cells.Item(6, 17).Value = 0.2115384615384615
$ws.Cells.Item(6, 18).Value = 0.07692307692307693
$ws.Cells.Item(6, 19).Value = 0.3413461538461539
$ws.Cells.Item(7, 2).Value = 0.120253164556962
$ws.Cells.Item(7, 4).Value = 0.0189873417721519
$ws.Cells.Item(7, 5).Value = 0.006329113924050633
$ws.Cells.Item(7, 6).Value = 0.05696202531645569
$ws.Cells.Item(7, 10).Value = 0.1329113924050633
$ws.Cells.Item(7, 15).Value = 0.02531645569620253
$ws.Cells.Item(7, 17).Value = 0.1518987341772152
$ws.Cells.Item(7, 18).Value = 0.08227848101265822
$ws.Cells.Item(7, 19).Value = 0.4050632911392405
$ws.Cells.Item(8, 2).Value = 0.08951965065502183
$ws.Cells.Item(8, 4).Value = 0.01091703056768559
$ws.Cells.Item(8, 5).Value = 0.004366812227074236
$ws.Cells.Item(8, 6).Value = 0.04366812227074236
$ws.Cells.Item(8, 10).Value = 0.1200873362445415
$ws.Cells.Item(8, 15).Value = 0.01310043668122271
$ws.Cells.Item(8, 17).Value = 0.1768558951965065
$ws.Cells.Item(8, 18).Value = 0.08951965065502183
$ws.Cells.Item(8, 19).Value = 0.4519650655021834
$ws.Cells.Item(9, 2).Value = 0.0860655737704918
$ws.Cells.Item(9, 4).Value = 0.01639344262295082
$ws.Cells.Item(9, 5).Value = 0.004098360655737705
$ws.Cells.Item(9, 6).Value = 0.05737704918032787
$ws.Cells.Item(9, 10).Value = 0.09836065573770492
$ws.Cells.Item(9, 15).Value = 0.02868852459016394
$ws.Cells.Item(9, 17).Value = 0.2336065573770492
$ws.Cells.Item(9, 18).Value = 0.04918032786885246
$ws.Cells.Item(9, 19).Value = 0.4262295081967213
$ws.Cells.Item(10, 2).Value = 0.1055853098699311
$ws.Cells.Item(10, 4).Value = 0.02065799540933435
$ws.Cells.Item(10, 5).Value = 0.0007651109410864575
$ws.Cells.Item(10, 6).Value = 0.06579954093343535
$ws.Cells.Item(10, 10).Value = 0.1185921958684009
$ws.Cells.Item(10, 15).Value = 0.01606732976281561
$ws.Cells.Item(10, 17).Value = 0.2471308339709258
$ws.Cells.Item(10, 18).Value = 0.08339709257842387
$ws.Cells.Item(10, 19).Value = 0.3420045906656465
$ws.Cells.Item(11, 6).Value = 0.004065040650406504
$ws.Cells.Item(11, 7).Value = 0.1382113821138211
$ws.Cells.Item(11, 10).Value = 0.09349593495934959
$ws.Cells.Item(11, 11).Value = 0.2195121951219512
$ws.Cells.Item(11, 12).Value = 0.5284552845528455
$ws.Cells.Item(11, 19).Value = 0.01626016260162602
$ws.Cells.Item(12, 7).Value = 0.7803030303030303
$ws.Cells.Item(12, 10).Value = 0.1742424242424243
$ws.Cells.Item(12, 12).Value = 0.02272727272727273
$ws.Cells.Item(12, 19).Value = 0.02272727272727273
$ws.Cells.Item(13, 7).Value = 0.6571428571428571
$ws.Cells.Item(13, 10).Value = 0.3142857142857143
$ws.Cells.Item(13, 19).Value = 0.02857142857142857
$ws.Cells.Item(15, 6).Value = 0.01574803149606299
$ws.Cells.Item(15, 8).Value = 0.1653543307086614
$ws.Cells.Item(15, 9).Value = 0.07874015748031496
$ws.Cells.Item(15, 10).Value = 0.3779527559055118
$ws.Cells.Item(15, 11).Value = 0.03937007874015748
$ws.Cells.Item(15, 13).Value = 0.01181102362204724
$ws.Cells.Item(15, 15).Value = 0.08267716535433071
$ws.Cells.Item(15, 19).Value = 0.2283464566929134
$ws.Cells.Item(16, 6).Value = 0.01578947368421053
$ws.Cells.Item(16, 8).Value = 0.1368421052631579
$ws.Cells.Item(16, 9).Value = 0.1421052631578947
$ws.Cells.Item(16, 10).Value = 0.4052631578947368
$ws.Cells.Item(16, 11).Value = 0.1105263157894737
$ws.Cells.Item(16, 13).Value = 0.02105263157894737
$ws.Cells.Item(16, 15).Value = 0.08421052631578947
$ws.Cells.Item(16, 19).Value = 0.08421052631578947
$ws.Cells.Item(17, 6).Value = 0.007561436672967864
$ws.Cells.Item(17, 8).Value = 0.1776937618147448
$ws.Cells.Item(17, 9).Value = 0.1077504725897921
$ws.Cells.Item(17, 10).Value = 0.44234404536862
$ws.Cells.Item(17, 11).Value = 0.08695652173913043
$ws.Cells.Item(17, 13).Value = 0.007561436672967864
$ws.Cells.Item(17, 14).Value = 0.005671077504725898
$ws.Cells.Item(17, 15).Value = 0.05671077504725898
$ws.Cells.Item(17, 19).Value = 0.1077504725897921
$ws.Cells.Item(18, 6).Value = 0.02631578947368421
$ws.Cells.Item(18, 8).Value = 0.1947368421052632
$ws.Cells.Item(18, 9).Value = 0.09473684210526316
$ws.Cells.Item(18, 10).Value = 0.4473684210526316
$ws.Cells.Item(18, 11).Value = 0.04736842105263158
$ws.Cells.Item(18, 13).Value = 0.03157894736842105
$ws.Cells.Item(18, 15).Value = 0.05789473684210526
$ws.Cells.Item(18, 19).Value = 0.1
$ws.Cells.Item(19, 6).Value = 0.0211038961038961
$ws.Cells.Item(19, 8).Value = 0.213474025974026
$ws.Cells.Item(19, 9).Value = 0.09983766233766234
$ws.Cells.Item(19, 10).Value = 0.3766233766233766
$ws.Cells.Item(19, 11).Value = 0.08847402597402597
$ws.Cells.Item(19, 13).Value = 0.01623376623376623
$ws.Cells.Item(19, 14).Value = 0.002435064935064935
$ws.Cells.Item(19, 15).Value = 0.08603896103896104
$ws.Cells.Item(19, 19).Value = 0.09577922077922078
